$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Add 39 new rows to the table (extends table ref automatically)
for ($i = 0; $i -lt 39; $i++) {
    $lo.ListRows.Add() | Out-Null
}

# Populate new rows with data
$ws.Range("A345").Value = "ALEXANDRE DA COSTA XAVIER"
$ws.Range("B345").Value = "SIL3D65"
$ws.Range("C345").Value = 2393826.0
$ws.Range("D345").Value = "SIL3D65"
$ws.Range("E345").Value = "Pedreira"

$ws.Range("A346").Value = "ANA ALICE VASCONCELOS DE JESUS"
$ws.Range("B346").Value = "JVY6337"
$ws.Range("C346").Value = 91931.0
$ws.Range("D346").Value = "JVY6337"
$ws.Range("E346").Value = "Mangueirão"

$ws.Range("A347").Value = "ANTONIO RONILSON NOGUEIRA REIS"
$ws.Range("B347").Value = "QDU5B51"
$ws.Range("C347").Value = 149588.0
$ws.Range("D347").Value = "QDU5B51"
$ws.Range("E347").Value = "Águas Lindas"

$ws.Range("A348").Value = "ARISON QUADROS DA CONCEICAO"
$ws.Range("B348").Value = "TCK6C61"
$ws.Range("C348").Value = 2387582.0
$ws.Range("D348").Value = "TCK6C61"
$ws.Range("E348").Value = "Parque Verde"

$ws.Range("A349").Value = "BRENA DE JESUS VIEIRA"
$ws.Range("B349").Value = "QEU2149"
$ws.Range("C349").Value = 387222.0
$ws.Range("D349").Value = "QEU2149"
$ws.Range("E349").Value = "Quarenta Horas (Coqueiro)"

$ws.Range("A350").Value = "CARLA ROBERTA RODRIGUES SEPEDA"
$ws.Range("B350").Value = "OFL8E96"
$ws.Range("C350").Value = 1185016.0
$ws.Range("D350").Value = "OFL8E96"
$ws.Range("E350").Value = "Cidade Nova"

$ws.Range("A351").Value = "CARLOS DANIELSON SANTOS DE ARAUJO"
$ws.Range("B351").Value = "QEO7216"
$ws.Range("C351").Value = 118538.0
$ws.Range("D351").Value = "QEO7216"
$ws.Range("E351").Value = "Maracangalha"

$ws.Range("A352").Value = "CARLOS MALCHER TEIXEIRA"
$ws.Range("B352").Value = "QZS5C57"
$ws.Range("C352").Value = 1270939.0
$ws.Range("D352").Value = "QZS5C57"
$ws.Range("E352").Value = "Cabanagem"

$ws.Range("A353").Value = "DANIEL DA COSTA BASTOS SILVA"
$ws.Range("B353").Value = "QEN7554"
$ws.Range("C353").Value = 1790487.0
$ws.Range("D353").Value = "QEN7554"
$ws.Range("E353").Value = "Santa Rita de Cássia"

$ws.Range("A354").Value = "DIEGO SALDANHA CORREIA"
$ws.Range("B354").Value = "QVD4083"
$ws.Range("C354").Value = 1804262.0
$ws.Range("D354").Value = "QVD4083"
$ws.Range("E354").Value = "Coqueiro"

$ws.Range("A355").Value = "DOUGLAS BARROSO DE ALMEIDA"
$ws.Range("B355").Value = "PBG7371"
$ws.Range("C355").Value = 2064106.0
$ws.Range("D355").Value = "PBG7371"
$ws.Range("E355").Value = "Marco"

$ws.Range("A356").Value = "ELAINE DO NASCIMENTO MACIEL"
$ws.Range("B356").Value = "QVH4455"
$ws.Range("C356").Value = 3022260.0
$ws.Range("D356").Value = "QVH4455"
$ws.Range("E356").Value = "Guamá"

$ws.Range("A357").Value = "ELIELSON DA COSTA SOARES"
$ws.Range("B357").Value = "QVT2H80"
$ws.Range("C357").Value = 91595.0
$ws.Range("D357").Value = "QVT2H80"
$ws.Range("E357").Value = "Icuí-Guajará"

$ws.Range("A358").Value = "EVERTON BARRETO BARBOSA"
$ws.Range("B358").Value = "NTA7E54"
$ws.Range("C358").Value = 421334.0
$ws.Range("D358").Value = "NTA7E54"
$ws.Range("E358").Value = "Carananduba (Mosqueiro)"

$ws.Range("A359").Value = "FABIO AUGUSTO SILVA"
$ws.Range("B359").Value = "QDK8792"
$ws.Range("C359").Value = 2899065.0
$ws.Range("D359").Value = "QDK8792"
$ws.Range("E359").Value = "Terra Firme"

$ws.Range("A360").Value = "FABIO XAVIER PORFIRIO"
$ws.Range("B360").Value = "Nes5b81"
$ws.Range("C360").Value = 817810.0
$ws.Range("D360").Value = "Nes5b81"
$ws.Range("E360").Value = "Marituba"

$ws.Range("A361").Value = "FAGNER KAIK SAMPAIO DOS SANTOS"
$ws.Range("B361").Value = "Jvw2867"
$ws.Range("C361").Value = 1057435.0
$ws.Range("D361").Value = "Jvw2867"
$ws.Range("E361").Value = "Marambaia"

$ws.Range("A362").Value = "GLESIELI CRISTINA DOS SANTOS"
$ws.Range("B362").Value = "QVJ8E67"
$ws.Range("C362").Value = 357281.0
$ws.Range("D362").Value = "QVJ8E67"
$ws.Range("E362").Value = "Guanabara"

$ws.Range("A363").Value = "GUSTAVO TEIXEIRA PEREIRA"
$ws.Range("B363").Value = "QVF4G65"
$ws.Range("C363").Value = 2449602.0
$ws.Range("D363").Value = "QVF4G65"
$ws.Range("E363").Value = "Tapanã (Icoaraci)"

$ws.Range("A364").Value = "HELIO FERNANDO DE LIMA MELO"
$ws.Range("B364").Value = "NSJ3B34"
$ws.Range("C364").Value = 724718.0
$ws.Range("D364").Value = "NSJ3B34"
$ws.Range("E364").Value = "Marambaia"

$ws.Range("A365").Value = "IVANIO SOARES MODESTO"
$ws.Range("B365").Value = "TDZ2E11"
$ws.Range("C365").Value = 2834712.0
$ws.Range("D365").Value = "TDZ2E11"
$ws.Range("E365").Value = "Maguari"

$ws.Range("A366").Value = "JANDRESSON GAIA DA SILVA"
$ws.Range("B366").Value = "OLL2B94"
$ws.Range("C366").Value = 1348538.0
$ws.Range("D366").Value = "OLL2B94"
$ws.Range("E366").Value = "Vila (Mosqueiro)"

$ws.Range("A367").Value = "JOAO PAULO MORAES ADRIANO"
$ws.Range("B367").Value = "QDU2F82"
$ws.Range("C367").Value = 1068245.0
$ws.Range("D367").Value = "QDU2F82"
$ws.Range("E367").Value = "São Brás"

$ws.Range("A368").Value = "JOSÉ RICARDO SOARES DOS SANTOS JUNIOR"
$ws.Range("B368").Value = "QUC2D96"
$ws.Range("C368").Value = 146997.0
$ws.Range("D368").Value = "QUC2D96"
$ws.Range("E368").Value = "Tenoné"

$ws.Range("A369").Value = "LENDREL MACIEL DE MELO"
$ws.Range("B369").Value = "OTU3793"
$ws.Range("C369").Value = 2908016.0
$ws.Range("D369").Value = "OTU3793"
$ws.Range("E369").Value = "Guamá"

$ws.Range("A370").Value = "LUAN DE JESUS BRITO SOUSA"
$ws.Range("B370").Value = "QVP3F31"
$ws.Range("C370").Value = 1735935.0
$ws.Range("D370").Value = "QVP3F31"
$ws.Range("E370").Value = "Sacramenta"

$ws.Range("A371").Value = "LUIZ FERNANDO DE ALMEIDA SANTOS"
$ws.Range("B371").Value = "OTN2394"
$ws.Range("C371").Value = 199477.0
$ws.Range("D371").Value = "OTN2394"
$ws.Range("E371").Value = "Coqueiro"

$ws.Range("A372").Value = "MATHEUS DE LUCAS SOUZA E SILVA"
$ws.Range("B372").Value = "OTA3F25"
$ws.Range("C372").Value = 406791.0
$ws.Range("D372").Value = "OTA3F25"
$ws.Range("E372").Value = "Centro"

$ws.Range("A373").Value = "MICILENE QUEIROZ QUARESMA"
$ws.Range("B373").Value = "TVN9J56"
$ws.Range("C373").Value = 132324.0
$ws.Range("D373").Value = "TVN9J56"
$ws.Range("E373").Value = "Telégrafo Sem Fio"

$ws.Range("A374").Value = "RICHARD CUNHA DA SILVA"
$ws.Range("B374").Value = "OSY3B52"
$ws.Range("C374").Value = 2412759.0
$ws.Range("D374").Value = "OSY3B52"
$ws.Range("E374").Value = "Cruzeiro (Icoaraci)"

$ws.Range("A375").Value = "RICHARD PATRICK PANTOJA COSTA"
$ws.Range("B375").Value = "QEO8E69"
$ws.Range("C375").Value = 201882.0
$ws.Range("D375").Value = "QEO8E69"
$ws.Range("E375").Value = "Coqueiro"

$ws.Range("A376").Value = "RINALDO DA ROCHA MENDES"
$ws.Range("B376").Value = "QVB4F50"
$ws.Range("C376").Value = 554581.0
$ws.Range("D376").Value = "QVB4F50"
$ws.Range("E376").Value = "Centro"

$ws.Range("A377").Value = "ROBERTO ALCIDES TELES LEAL"
$ws.Range("B377").Value = "QEI9432"
$ws.Range("C377").Value = 92308.0
$ws.Range("D377").Value = "QEI9432"
$ws.Range("E377").Value = "Umarizal"

$ws.Range("A378").Value = "ROMULO HENRIQUE ARAUJO GONCALVES"
$ws.Range("B378").Value = "NSF4F52"
$ws.Range("C378").Value = 1432400.0
$ws.Range("D378").Value = "NSF4F52"
$ws.Range("E378").Value = "Atalaia"

$ws.Range("A379").Value = "RUAN VITOR LIMA SILVA"
$ws.Range("B379").Value = "RUU5J86"
$ws.Range("C379").Value = 1750026.0
$ws.Range("D379").Value = "RUU5J86"
$ws.Range("E379").Value = "Nazaré"

$ws.Range("A380").Value = "SILAS PATRICK PEREIRA DE OLIVEIRA"
$ws.Range("B380").Value = "JVD7661"
$ws.Range("C380").Value = 1165098.0
$ws.Range("D380").Value = "JVD7661"
$ws.Range("E380").Value = "Bengui"

$ws.Range("A381").Value = "VICTOR MANOEL DOS SANTOS MATOS"
$ws.Range("B381").Value = "Rvw6e87"
$ws.Range("C381").Value = 584397.0
$ws.Range("D381").Value = "Rvw6e87"
$ws.Range("E381").Value = "Castanheira"

$ws.Range("A382").Value = "WILLIAMS TIAGO FARIAS PINHEIRO"
$ws.Range("B382").Value = "QEE7B13"
$ws.Range("C382").Value = 1792651.0
$ws.Range("D382").Value = "QEE7B13"
$ws.Range("E382").Value = "Pedreira"

$ws.Range("A383").Value = "YORDANIS GONZALEZ CRUZ"
$ws.Range("B383").Value = "RNG3J95"
$ws.Range("C383").Value = 264574.0
$ws.Range("D383").Value = "RNG3J95"
$ws.Range("E383").Value = "Icuí-Guajará"

# Copy formatting from the last original data row (344) to all new rows (345:383)
$ws.Range("A344:E344").Copy()
$ws.Range("A345:E383").PasteSpecial(-4122)
$excel.CutCopyMode = 0

Write-Host "Table range now: $($lo.Range.Address())"